$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openTickets")

# Update the highlighting-annotation JSON text in column G (ticketDescriptionHighlighting)
# for rows 2-4. The text is reformatted (no space after the JSON ":" separators) and,
# for row 3, a stray key/value separator is dropped ("key" "Service Anfrage").
$ws.Range("G2").Value = '[{"start":143, "end":147, "key":"System"}​,{ "start":104, "end":128, "key":"Fehlerbeschreibung"}​,{ "start":67,"end":77, "key":"System"}]'
$ws.Range("G3").Value = '[{"start":229, "end":297 ,"key" "Service Anfrage"},​{ "start":191, "end":192, "key":"System"}​,{"start":176, "end":191, "key":"System"}​,{"start":129, "end":144, "key":"System"}]'
$ws.Range("G4").Value = '[{"start":130, "end":165, "key":"Auslöser"}​,{ "start":37, "end":78, "key":"Fehlerbeschreibung"},​{ "start":24, "end":36, "key":"System"}]'

# Move the active selection from G5 to G4.
$ws.Range("G4").Select()
